$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D4 and D5 cluster labels (Inflammatory-Mac row removed; MuSCs/Resolving-Mac shift up)
$ws.Range("D4").Value = "MuSCs"
$ws.Range("D5").Value = "Resolving-Mac"

# Update recalculated TPM-derived values

# Row 2
$ws.Range("M2").Value = 0.08425766666666668
$ws.Range("N2").Value = 0.252773
$ws.Range("O2").Value = 0.007654801123801229
$ws.Range("P2").Value = 0.008027752567511702
$ws.Range("Q2").Value = 0.0327257900768889
$ws.Range("R2").Value = 0.2945321106920001
$ws.Range("S2").Value = 0.007654801123801229
$ws.Range("T2").Value = 0.008027752567511702

# Row 3
$ws.Range("O3").Value = 0.8522349591772004
$ws.Range("P3").Value = 0.8937569077249424
$ws.Range("S3").Value = 0.8522349591772004
$ws.Range("T3").Value = 0.8937569077249424

# Row 4
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.534105
$ws.Range("N4").Value = 3.06821
$ws.Range("O4").Value = 0.139373295542195
$ws.Range("P4").Value = 0.09744249071366434
$ws.Range("Q4").Value = 0.5958484274733334
$ws.Range("R4").Value = 3.57509056484
$ws.Range("S4").Value = 0.139373295542195
$ws.Range("T4").Value = 0.09744249071366434

# Row 5
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.008111666666666666
$ws.Range("N5").Value = 0.024335
$ws.Range("O5").Value = 0.0007369441568035466
$ws.Range("P5").Value = 0.0007728489938814559
$ws.Range("Q5").Value = 0.003150582148888889
$ws.Range("R5").Value = 0.02835523934
$ws.Range("S5").Value = 0.0007369441568035466
$ws.Range("T5").Value = 0.0007728489938814559

# Remove the now-obsolete Inflammatory-Mac row (old row 6)
$ws.Rows("6:6").Delete()
